$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (I1, J1)
$ws.Range("I1").Value = "Other found locations"
$ws.Range("J1").Value = "Misc. Data"

# Row 2: update Authors (E2) with refreshed list, and populate new columns
$ws.Range("E2").Value = @'
[Christiane Maria Prado%Jeronimo%NULL%0,   Maria Eduarda Leão%Farias%NULL%1,   Fernando Fonseca Almeida%Val%NULL%1,   Vanderson Souza%Sampaio%NULL%1,   Marcia Almeida Araújo%Alexandre%NULL%1,   Gisely Cardoso%Melo%NULL%1,   Izabella Picinin%Safe%NULL%1,   Mayla Gabriela Silva%Borba%NULL%1,   Rebeca Linhares%Abreu-Netto%NULL%1,   Alex Bezerra Silva%Maciel%NULL%1,   João Ricardo Silva%Neto%NULL%1,   Lucas Barbosa%Oliveira%NULL%1,   Erick Frota Gomes%Figueiredo%NULL%1,   Kelry Mazurega Oliveira%Dinelly%NULL%1,   Maria Gabriela de Almeida%Rodrigues%NULL%1,   Marcelo%Brito%NULL%1,   Maria Paula Gomes%Mourão%NULL%1,   Guilherme Augusto%Pivoto João%NULL%1,   Ludhmila Abrahão%Hajjar%NULL%1,   Quique%Bassat%NULL%1,   Gustavo Adolfo Sierra%Romero%NULL%1,   Felipe Gomes%Naveca%NULL%1,   Heline Lira%Vasconcelos%NULL%1,   Michel de Araújo%Tavares%NULL%1,   José Diego%Brito-Sousa%NULL%1,   Fabio Trindade Maranhão%Costa%NULL%1,   Maurício Lacerda%Nogueira%NULL%1,   Djane%Baía-da-Silva%NULL%1,   Mariana Simão%Xavier%NULL%1,   Wuelton Marcelo%Monteiro%NULL%1,   Marcus Vinícius Guimarães%Lacerda%marcuslacerda.br@gmail.com%1,   NULL%NULL%NULL%1]
'@
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""

# Row 3: only populate new columns (Authors E3 unchanged)
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""

# Row 4: update Authors (E4) with refreshed list, and populate new columns
$ws.Range("E4").Value = @'
[Alexandre%B. Cavalcanti%coreGivesNoEmail%0,  Alvaro%Avezum%coreGivesNoEmail%0,  Andre%N. Costa%coreGivesNoEmail%0,  Bruno%M. Tomazini%coreGivesNoEmail%0,  Caio%C. F. Fernandes%coreGivesNoEmail%0,  Cassia%Righy%coreGivesNoEmail%0,  Cristina%P. Amendola%coreGivesNoEmail%0,  Daniel%N. Forte%coreGivesNoEmail%0,  Daniela%H. M. Freitas%coreGivesNoEmail%0,  Douglas%Costa Morais%coreGivesNoEmail%0,  Eduardo%L. V. Costa%coreGivesNoEmail%0,  Fernando%G. Zampieri%coreGivesNoEmail%0,  Flavia%R. Bueno%coreGivesNoEmail%0,  Fl\u00e1via%R. Machado%coreGivesNoEmail%0,  Fl\u00e1vio%G. R. Freitas%coreGivesNoEmail%0,  Franca%P. Baldassare%coreGivesNoEmail%0,  Gedealvares%F. S. Junior%coreGivesNoEmail%0,  Guilherme%B. Olivato%coreGivesNoEmail%0,  Israel%S. Maia%coreGivesNoEmail%0,  Let\u00edcia%Kawano-Dourado%coreGivesNoEmail%0,  Livia%M. G. Melro%coreGivesNoEmail%0,  Lucas%P. Damiani%coreGivesNoEmail%0,  Luciano%C. P. Azevedo%coreGivesNoEmail%0,  Maria%Vitoria A. O. Silva%coreGivesNoEmail%0,  Michele%O. Honorato%coreGivesNoEmail%0,  Otavio%Berwanger%coreGivesNoEmail%0,  Regis%G. Rosa%coreGivesNoEmail%0,  Renato%D. Lopes%coreGivesNoEmail%0,  Ricardo%A. B. Moura%coreGivesNoEmail%0,  Roberta%M. L. Roepke%coreGivesNoEmail%0,  Stevin%Zung%coreGivesNoEmail%0,  Thiago%Lisboa%coreGivesNoEmail%0,  Viviane%C. Veiga%coreGivesNoEmail%0]
'@
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""

# Row 5: Authors (E5) field was too large and is now cleared entirely, populate new columns
$ws.Range("E5").ClearContents()
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""

# Row 6: update Authors (E6) with refreshed list, and populate new columns
$ws.Range("E6").Value = @'
[Annane%Djillali%coreGivesNoEmail%0,  Aubron%C\u00e9cile%coreGivesNoEmail%0,  Badi\u00e9%Julio%coreGivesNoEmail%0,  Bourgoin%H\u00e9l\u00e8ne%coreGivesNoEmail%0,  Caille-F\u00e9n\u00e9rol%Caroline%coreGivesNoEmail%0,  Coffre%Carine%coreGivesNoEmail%0,  Dequin%Pierre-Fran\u00e7ois%coreGivesNoEmail%0,  Ehrmann%Stephan%coreGivesNoEmail%0,  Fran\u00e7ois%Bruno%coreGivesNoEmail%0,  Giraudeau%Bruno%coreGivesNoEmail%0,  Guillon%Antoine%coreGivesNoEmail%0,  Heming%Nicholas%coreGivesNoEmail%0,  Jouan%Youenn%coreGivesNoEmail%0,  Le%Gouge Am\u00e9lie%coreGivesNoEmail%0,  Leclerc%Marie%coreGivesNoEmail%0,  Lengell\u00e9%C\u00e9line%coreGivesNoEmail%0,  Meziani%Ferhat%coreGivesNoEmail%0,  Plantef\u00e8ve%Ga\u00ebtan%coreGivesNoEmail%0,  Ricard%Jean-Damien%coreGivesNoEmail%0,  Tavernier%Elsa%coreGivesNoEmail%0,  Voiriot%Guillaume%coreGivesNoEmail%0,  Zohar%Sarah%coreGivesNoEmail%0]
'@
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
